$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 3 with the new ticket entry
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "ClassPlus"
$ws.Range("C3").Value = 1285491

# Copy the date formatting from D2 so D3 reuses the same (date) style,
# then overwrite with the actual serial date value.
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 45208

$ws.Range("E3").Value = "Anirban Chakraborty"
$ws.Range("F3").Value = " Debasish Nath cannot Login as Sub Admin"
$ws.Range("G3").Value = "Pending"

# Update the active selection to match the edited document
$ws.Range("F7").Select()
